{"js": "// Update the date and the ten division problems for this worksheet.\nconst replacements = [\n  [\"2023-09-18 Monday\", \"2023-09-19 Tuesday\"],\n  [\"96\u00f73=\", \"78\u00f76=\"],\n  [\"62\u00f79=\", \"20\u00f73=\"],\n  [\"83\u00f73=\", \"20\u00f77=\"],\n  [\"99\u00f73=\", \"50\u00f73=\"],\n  [\"76\u00f73=\", \"10\u00f77=\"],\n  [\"57\u00f75=\", \"38\u00f72=\"],\n  [\"88\u00f75=\", \"40\u00f75=\"],\n  [\"92\u00f75=\", \"98\u00f79=\"],\n  [\"46\u00f78=\", \"76\u00f79=\"],\n  [\"54\u00f72=\", \"27\u00f76=\"],\n  [\"33\u00f74=\", \"71\u00f75=\"],\n  [\"76\u00f74=\", \"73\u00f73=\"],\n  [\"34\u00f77=\", \"78\u00f79=\"],\n  [\"98\u00f75=\", \"62\u00f78=\"],\n  [\"23\u00f79=\", \"24\u00f74=\"],\n  [\"20\u00f78=\", \"10\u00f75=\"],\n  [\"18\u00f77=\", \"51\u00f78=\"],\n  [\"33\u00f75=\", \"72\u00f78=\"],\n  [\"98\u00f78=\", \"44\u00f76=\"],\n  [\"34\u00f73=\", \"65\u00f72=\"],\n  [\"11\u00f77=\", \"35\u00f77=\"],\n  [\"19\u00f72=\", \"95\u00f74=\"],\n  [\"87\u00f77=\", \"29\u00f73=\"],\n  [\"58\u00f78=\", \"10\u00f73=\"],\n  [\"37\u00f72=\", \"75\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and the ten division problems for this worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-18 Monday\", \"2023-09-19 Tuesday\"),\n    @(\"96\u00f73=\", \"78\u00f76=\"),\n    @(\"62\u00f79=\", \"20\u00f73=\"),\n    @(\"83\u00f73=\", \"20\u00f77=\"),\n    @(\"99\u00f73=\", \"50\u00f73=\"),\n    @(\"76\u00f73=\", \"10\u00f77=\"),\n    @(\"57\u00f75=\", \"38\u00f72=\"),\n    @(\"88\u00f75=\", \"40\u00f75=\"),\n    @(\"92\u00f75=\", \"98\u00f79=\"),\n    @(\"46\u00f78=\", \"76\u00f79=\"),\n    @(\"54\u00f72=\", \"27\u00f76=\"),\n    @(\"33\u00f74=\", \"71\u00f75=\"),\n    @(\"76\u00f74=\", \"73\u00f73=\"),\n    @(\"34\u00f77=\", \"78\u00f79=\"),\n    @(\"98\u00f75=\", \"62\u00f78=\"),\n    @(\"23\u00f79=\", \"24\u00f74=\"),\n    @(\"20\u00f78=\", \"10\u00f75=\"),\n    @(\"18\u00f77=\", \"51\u00f78=\"),\n    @(\"33\u00f75=\", \"72\u00f78=\"),\n    @(\"98\u00f78=\", \"44\u00f76=\"),\n    @(\"34\u00f73=\", \"65\u00f72=\"),\n    @(\"11\u00f77=\", \"35\u00f77=\"),\n    @(\"19\u00f72=\", \"95\u00f74=\"),\n    @(\"87\u00f77=\", \"29\u00f73=\"),\n    @(\"58\u00f78=\", \"10\u00f73=\"),\n    @(\"37\u00f72=\", \"75\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$findText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$true, [ref]$replaceText, [ref]2)\n}\n"}
